$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.230.08"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.027.24"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.34"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.33"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0791"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  -4.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.326.64"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.24"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.40"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.17"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.026.70"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.178.02"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.30"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.24"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.61"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  -5.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.30"
$ws.Range("E26").Value = "  -6.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.52"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.127"
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.36"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.82"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.53"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.46"
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.49"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.479.98"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.28"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0922"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.39"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.02"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.24"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.212.75"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.62"
$ws.Range("E51").Value = "  -9.90%  "
